$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# M2 is an empty styled cell (matches L2's border-only style)
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial($xlFormats)

$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial($xlFormats)
$ws.Range("M3").Value = 2021

$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial($xlFormats)
$ws.Range("M4").Value = 2.0173148373954581

$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial($xlFormats)
$ws.Range("M5").Value = 0.11867182493532386

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial($xlFormats)
$ws.Range("M6").Value = 3.9440914499323179

$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial($xlFormats)
$ws.Range("M7").Value = 0

$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial($xlFormats)
$ws.Range("M8").Value = "-"

$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial($xlFormats)
$ws.Range("M9").Value = 0

$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial($xlFormats)
$ws.Range("M10").Value = 0.62921030174566528

$ws.Range("L11").Copy()
$ws.Range("M11").PasteSpecial($xlFormats)
$ws.Range("M11").Value = "-"

$ws.Range("L12").Copy()
$ws.Range("M12").PasteSpecial($xlFormats)
$ws.Range("M12").Value = 1.2497227177719943

$ws.Range("L13").Copy()
$ws.Range("M13").PasteSpecial($xlFormats)
$ws.Range("M13").Value = 0.19844537890168421

$ws.Range("L14").Copy()
$ws.Range("M14").PasteSpecial($xlFormats)
$ws.Range("M14").Value = "-"

$ws.Range("L15").Copy()
$ws.Range("M15").PasteSpecial($xlFormats)
$ws.Range("M15").Value = 0.39861918314956984

$ws.Range("L16").Copy()
$ws.Range("M16").PasteSpecial($xlFormats)
$ws.Range("M16").Value = 0

$ws.Range("L17").Copy()
$ws.Range("M17").PasteSpecial($xlFormats)
$ws.Range("M17").Value = "-"

$ws.Range("L18").Copy()
$ws.Range("M18").PasteSpecial($xlFormats)
$ws.Range("M18").Value = 0

$ws.Range("L19").Copy()
$ws.Range("M19").PasteSpecial($xlFormats)
$ws.Range("M19").Value = 0.85521252031129735

$ws.Range("L20").Copy()
$ws.Range("M20").PasteSpecial($xlFormats)
$ws.Range("M20").Value = "-"

$ws.Range("L21").Copy()
$ws.Range("M21").PasteSpecial($xlFormats)
$ws.Range("M21").Value = 1.6913581464969858

$ws.Range("L22").Copy()
$ws.Range("M22").PasteSpecial($xlFormats)
$ws.Range("M22").Value = 1.8347815875998121

$ws.Range("L23").Copy()
$ws.Range("M23").PasteSpecial($xlFormats)
$ws.Range("M23").Value = "-"

$ws.Range("L24").Copy()
$ws.Range("M24").PasteSpecial($xlFormats)
$ws.Range("M24").Value = 3.6321107648498847

$ws.Range("L25").Copy()
$ws.Range("M25").PasteSpecial($xlFormats)
$ws.Range("M25").Value = 6.1211560415300026

$ws.Range("L26").Copy()
$ws.Range("M26").PasteSpecial($xlFormats)
$ws.Range("M26").Value = "-"
$ws.Range("M26").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

$ws.Range("L27").Copy()
$ws.Range("M27").PasteSpecial($xlFormats)
$ws.Range("M27").Value = 12.437939862560766

$ws.Range("L28").Copy()
$ws.Range("M28").PasteSpecial($xlFormats)
$ws.Range("M28").Value = 3.6823562661275693

$ws.Range("L29").Copy()
$ws.Range("M29").PasteSpecial($xlFormats)
$ws.Range("M29").Value = 0.69433233870225819

$ws.Range("L30").Copy()
$ws.Range("M30").PasteSpecial($xlFormats)
$ws.Range("M30").Value = 7.0564990356117976

$ws.Range("L31").Copy()
$ws.Range("M31").PasteSpecial($xlFormats)
$ws.Range("M31").Value = 2.7447727328177227

$ws.Range("L8").Copy()
$ws.Range("M32").PasteSpecial($xlFormats)
$ws.Range("M32").Value = "-"

$ws.Range("L33").Copy()
$ws.Range("M33").PasteSpecial($xlFormats)
$ws.Range("M33").Value = 5.6418550419377889

$ws.Range("P6").Select()
$excel.CutCopyMode = $false
